$p = $ppt.ActivePresentation

# Insert a new slide at position 15 (before the current "Regression Analysis"
# slide), using the same "Title and Content" layout (layout index 2) as its
# neighbours so it gets a Title placeholder + a Content placeholder.
$newSlide = $p.Slides.Add(15, 2)

# --- Title placeholder -------------------------------------------------
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Hypothesis Testing"
$titleRange.LanguageID = "en-US"
$titleRange.ParagraphFormat.Alignment = 2

# --- Content placeholder ------------------------------------------------
$bodyShape = $newSlide.Shapes.Item(2)
$bodyShape.TextFrame.AutoSize = 2

$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "Conducted hypothesis testing to determine if diabetes contributes to death by heart failure`rNull Hypothesis: Diabetes does not contribute to death by heart failure `rP Value is 0.926 which indicates that the diabetes variable is not statistically significant, since P value is not less then 0.05`rThis validates the Null Hypothesis that diabetes does not contribute to death by heart failure "
$bodyRange.LanguageID = "en-US"
